# Add two new columns ("Upload" / "Status") after the existing "File"
# column, relabel the first header ("Label" -> "Recipe"), and mark every
# existing data row as already uploaded ("Y") leaving Status blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old Name/ID/YYYY/MM/DD/Label columns (C:G) two to the right,
# opening up C:D for the new "Upload"/"Status" columns.
$ws.Columns("C:D").Insert()

# Header row
$ws.Range("A1").Value = "Recipe"
$ws.Range("C1").Value = "Upload"
$ws.Range("D1").Value = "Status"

# Mark all existing example rows as already uploaded; Status stays blank.
$ws.Range("C2:C7").Value = "Y"
